$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.191741704940796
$ws.Range("B1").Value = 4.286515235900879
$ws.Range("C1").Value = 2.248613834381104
$ws.Range("D1").Value = 1.741935729980469
$ws.Range("E1").Value = 1.578306555747986
